$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 100.333336
$ws.Range("I2").Value = 100.333336
$ws.Range("K2").Value = 100.333336
$ws.Range("M2").Value = 12.666664

# row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()

# row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# row 29
$ws.Range("H29").Value = 3785.1428
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 3785.1428
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 11355.4284
$ws.Range("N29").Value = -11917.4284
$ws.Range("M29").ClearContents()

# row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# row 62
$ws.Range("H62").Value = 2346.5625
$ws.Range("I62").Value = 2159.8
$ws.Range("J62").Value = 2657.8333
$ws.Range("K62").Value = 2159.8
$ws.Range("L62").Value = 2657.8333
$ws.Range("M62").Value = -1535.8
$ws.Range("N62").Value = -3905.8333

# row 65
$ws.Range("H65").Value = 2346.5625
$ws.Range("I65").Value = 2159.8
$ws.Range("J65").Value = 2657.8333
$ws.Range("K65").Value = 10799
$ws.Range("L65").Value = 13289.1665
$ws.Range("M65").Value = -7679
$ws.Range("N65").Value = -19529.1665

# row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# row 110
$ws.Range("H110").Value = 45725
$ws.Range("J110").Value = 45725
$ws.Range("L110").Value = 45725
$ws.Range("N110").Value = -53905

# row 129
$ws.Range("H129").Value = 814.1
$ws.Range("J129").Value = 899
$ws.Range("L129").Value = 2697
$ws.Range("N129").Value = -12697

# row 137
$ws.Range("H137").Value = 32034.559
$ws.Range("I137").Value = 2965.111
$ws.Range("J137").Value = 64737.688
$ws.Range("K137").Value = 8895.332999999999
$ws.Range("L137").Value = 194213.064
$ws.Range("M137").Value = -6345.332999999999
$ws.Range("N137").Value = -199313.064

# row 138
$ws.Range("H138").Value = 3074.5278
$ws.Range("I138").Value = 1524.25
$ws.Range("J138").Value = 3268.3125
$ws.Range("K138").Value = 4572.75
$ws.Range("L138").Value = 9804.9375
$ws.Range("M138").Value = 567.25
$ws.Range("N138").Value = -20084.9375

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 18619.588
$ws.Range("I32").Value = 22902.416
$ws.Range("J32").Value = 4914.533
$ws.Range("K32").Value = 22902.416
$ws.Range("L32").Value = 4914.533
$ws.Range("M32").Value = -22615.416
$ws.Range("N32").Value = -5488.533

# row 61
$ws.Range("H61").Value = 377921.6
$ws.Range("I61").Value = 564894.4399999999
$ws.Range("J61").Value = 3975.875
$ws.Range("K61").Value = 564894.4399999999
$ws.Range("L61").Value = 3975.875
$ws.Range("M61").Value = -564682.4399999999
$ws.Range("N61").Value = -4399.875

# row 80
$ws.Range("H80").Value = 43440
$ws.Range("J80").Value = 43440
$ws.Range("L80").Value = 43440
$ws.Range("N80").Value = -45436

# row 83
$ws.Range("H83").Value = 43440
$ws.Range("J83").Value = 43440
$ws.Range("L83").Value = 130320
$ws.Range("N83").Value = -140304

# row 132
$ws.Range("H132").Value = 11623.49
$ws.Range("I132").Value = 1878.4048
$ws.Range("K132").Value = 5635.2144
$ws.Range("M132").Value = -3105.2144

# row 136
$ws.Range("H136").Value = 377921.6
$ws.Range("I136").Value = 564894.4399999999
$ws.Range("J136").Value = 3975.875
$ws.Range("K136").Value = 1694683.32
$ws.Range("L136").Value = 11927.625
$ws.Range("M136").Value = -1692133.32
$ws.Range("N136").Value = -17027.625

# row 139
$ws.Range("H139").Value = 51238.332
$ws.Range("J139").Value = 51238.332
$ws.Range("L139").Value = 51238.332
$ws.Range("N139").Value = -61518.332

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 1758.1818
$ws.Range("I99").Value = 1243.3334
$ws.Range("J99").Value = 2376
$ws.Range("K99").Value = 1243.3334
$ws.Range("L99").Value = 2376
$ws.Range("M99").Value = 254.6666
$ws.Range("N99").Value = -5372

# row 107
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# row 132
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120

# row 140
$ws.Range("H140").Value = 47685
$ws.Range("J140").Value = 47685
$ws.Range("L140").Value = 47685
$ws.Range("N140").Value = -58045

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1173
$ws.Range("I16").Value = 1020.2
$ws.Range("J16").Value = 1257.8889
$ws.Range("K16").Value = 1020.2
$ws.Range("L16").Value = 1257.8889
$ws.Range("M16").Value = -733.2
$ws.Range("N16").Value = -1831.8889

# row 31
$ws.Range("H31").Value = 12184.909
$ws.Range("I31").Value = 18683.084
$ws.Range("J31").Value = 4387.1
$ws.Range("K31").Value = 18683.084
$ws.Range("L31").Value = 4387.1
$ws.Range("M31").Value = -18388.084
$ws.Range("N31").Value = -4977.1

# row 34
$ws.Range("H34").Value = 12184.909
$ws.Range("I34").Value = 18683.084
$ws.Range("J34").Value = 4387.1
$ws.Range("K34").Value = 18683.084
$ws.Range("L34").Value = 4387.1
$ws.Range("M34").Value = -18481.084
$ws.Range("N34").Value = -4791.1

# row 52
$ws.Range("H52").Value = 38749.5
$ws.Range("J52").Value = 38749.5
$ws.Range("L52").Value = 38749.5
$ws.Range("N52").Value = -39337.5

# row 113
$ws.Range("H113").Value = 1173
$ws.Range("I113").Value = 1020.2
$ws.Range("J113").Value = 1257.8889
$ws.Range("K113").Value = 1020.2
$ws.Range("L113").Value = 1257.8889
$ws.Range("M113").Value = 1149.8
$ws.Range("N113").Value = -5597.8889

# row 132
$ws.Range("H132").Value = 19960.912
$ws.Range("I132").Value = 20827.709
$ws.Range("K132").Value = 62483.12699999999
$ws.Range("M132").Value = -59953.12699999999

$ws = $wb.Worksheets.Item("CUL")
# row 7
$ws.Range("H7").Value = 250.5
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 109
$ws.Range("N7").Value = -1724

# row 56
$ws.Range("H56").Value = 6550
$ws.Range("I56").Value = 6550
$ws.Range("K56").Value = 6550
$ws.Range("M56").Value = -6020

# row 68
$ws.Range("H68").Value = 4332.0938
$ws.Range("J68").Value = 4833.143
$ws.Range("L68").Value = 14499.429
$ws.Range("N68").Value = -16121.429

# row 71
$ws.Range("H71").Value = 4332.0938
$ws.Range("J71").Value = 4833.143
$ws.Range("L71").Value = 43498.287
$ws.Range("N71").Value = -51610.287

# row 92
$ws.Range("H92").Value = 62500350
$ws.Range("I92").Value = 62500350
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 187501050
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -187499802
$ws.Range("N92").ClearContents()

# row 107
$ws.Range("H107").Value = 5144.7915
$ws.Range("J107").Value = 1113.75
$ws.Range("L107").Value = 3341.25
$ws.Range("N107").Value = -7181.25

# row 131
$ws.Range("H131").Value = 173240.86
$ws.Range("J131").Value = 200839.6
$ws.Range("L131").Value = 602518.8
$ws.Range("N131").Value = -612598.8

# row 140
$ws.Range("H140").Value = 2310.5
$ws.Range("J140").Value = 3892.2
$ws.Range("L140").Value = 11676.6
$ws.Range("N140").Value = -22036.6

$ws = $wb.Worksheets.Item("GSM")
# row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 2072
$ws.Range("I22").Value = 2310
$ws.Range("K22").Value = 2310
$ws.Range("M22").Value = -2015

# row 27
$ws.Range("H27").Value = 2072
$ws.Range("I27").Value = 2310
$ws.Range("K27").Value = 2310
$ws.Range("M27").Value = -2203

# row 122
$ws.Range("H122").Value = 3032.2632
$ws.Range("I122").Value = 2551
$ws.Range("J122").Value = 3382.2727
$ws.Range("K122").Value = 7653
$ws.Range("L122").Value = 10146.8181
$ws.Range("M122").Value = -5203
$ws.Range("N122").Value = -15046.8181

# row 136
$ws.Range("H136").Value = 2006.5122
$ws.Range("I136").Value = 1479.52
$ws.Range("J136").Value = 2829.9375
$ws.Range("K136").Value = 4438.559999999999
$ws.Range("L136").Value = 8489.8125
$ws.Range("M136").Value = -1888.559999999999
$ws.Range("N136").Value = -13589.8125

$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("H136").Value = 1687.4762
$ws.Range("I136").Value = 1070.3334
$ws.Range("J136").Value = 2150.3333
$ws.Range("K136").Value = 3211.0002
$ws.Range("L136").Value = 6450.999899999999
$ws.Range("M136").Value = -661.0001999999999
$ws.Range("N136").Value = -11550.9999
